$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns AQ (43) and AR (44) in row 1
$ws.Cells.Item(1, 43).Value = "Pt"
$ws.Cells.Item(1, 44).Value = "Pt.Se"

# Copy the style of AP1 (col 42) onto the new header cells so they match the existing
# numeric-header formatting (style index 2 in the original workbook).
$ws.Cells.Item(1, 42).Copy() | Out-Null
$ws.Cells.Item(1, 43).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(1, 44).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data for AQ2:AR15
$ptValues = @(8.91, 8.5399999999999991, 7.72, 9.4700000000000006, 8.73, 9.74, 9.0399999999999991, 9.9, 9.09, 10.8, 8.8699999999999992, 9.89, 9.85, 9)
$ptSeValues = @(0.33, 0.44, 0.3, 0.39, 0.33, 0.31, 0.32, 0.55000000000000004, 0.26, 0.59, 0.6, 1.1200000000000001, 0.28000000000000003, 0.52)

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 43).Value = $ptValues[$i]
    $ws.Cells.Item($row, 44).Value = $ptSeValues[$i]
}

# Copy style from AP column data cells (style index 2) to new AQ/AR data cells
$ws.Cells.Item(2, 42).Copy() | Out-Null
$srcRange = $ws.Range($ws.Cells.Item(2, 43), $ws.Cells.Item(15, 44))
$srcRange.PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Change style of E2:F15 (style 7 -> style 1): remove the fill that used to be applied
$rangeEF = $ws.Range("E2:F15")
$rangeEF.Interior.Pattern = -4142  # xlPatternNone

# Update the sheet view (topLeftCell + selection) as reflected in the diff
$ws.Application.ActiveWindow.ScrollColumn = 28  # AB
$ws.Range("O23").Select()
